$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The GrantsTab query (same Cypher that used to back "ProjectsTab") now also
# reports a distinct "Grants" count in addition to "Projects".
$grantsTabQuery = @'
MATCH (pr:project)-->(pgm:program)
where pr.lead_doc='CCG'
WITH DISTINCT pr, pgm
RETURN
coalesce(pr.project_id , '') AS `Grant ID`,
coalesce(pr.queried_project_id, '') AS `Project ID`,
coalesce (pgm.program_id, '')AS `Program`,
coalesce(pr.project_title, '') AS `Project Title`,
coalesce(pr.principal_investigators, '') AS `Principal Investigators`,
coalesce(pr.program_officers, '') AS `Program Officers`,
coalesce(pr.lead_doc, '')AS `Lead DOC`,
SUBSTRING(pr.project_id, 1, 3) AS `Activity code`,
"$" + apoc.number.format(toInteger(pr.award_amount)) AS `Award Amount`,
coalesce(pr.project_end_date, '') AS `Project End Date`,
coalesce(pr.fiscal_year,'')AS `Fiscal Year`
'@

# New StatQuery (column C) shared by every row - replaces the old ProjectsTab
# StatQuery with one that separates Projects vs Grants and widens the dataset
# match to a variable-length hop.
$newStatQuery = @'
MATCH (p:program)<--(pr:project)
where pr.lead_doc='CCG'
OPTIONAL MATCH (pr)<--(pub:publication)
OPTIONAL MATCH (ct:clinical_trial)
WHERE EXISTS((pr)<--(pub)<--(ct)) OR EXISTS((pr)<--(ct))
OPTIONAL MATCH (pr)<--(pat)
WHERE pat:patent_application OR pat:granted_patent
OPTIONAL MATCH (pr)<-[*1..2]-(dt)
WHERE dt:sra OR dt:dbgap OR dt:geo
WITH p, pr, pub, ct, pat, dt
RETURN
COUNT(DISTINCT p.program_id) AS Programs,
COUNT(DISTINCT pr.queried_project_id) AS Projects,
COUNT(DISTINCT pr.project_id) AS Grants,
COUNT(DISTINCT pub.publication_id) AS Publications,
COUNT(DISTINCT dt.accession) AS Datasets,
COUNT(DISTINCT ct.clinical_trial_id) AS `Clinical Trials`,
COUNT(DISTINCT pat.patent_id) AS Patents
'@

# Row 2: ProjectsTab -> GrantsTab (the query in column B already matched this
# content; only the tab name and the shared StatQuery in column C change).
$ws.Range("A2").Value = "GrantsTab"
$ws.Range("B2").Value = $grantsTabQuery
$ws.Range("C2").Value = $newStatQuery

# Rows 3-6 keep their tab name (A) and dbExcel query (B); only the StatQuery
# in column C is refreshed to the new combined query.
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery
$ws.Range("C5").Value = $newStatQuery
$ws.Range("C6").Value = $newStatQuery

# Row 2 grew a line taller once wrapped (255 -> 270 points); other rows are
# unchanged.
$ws.Rows.Item(2).RowHeight = 270

# Selection/scroll moved to A6.
$ws.Range("A6").Select()
$excel.ActiveWindow.ScrollRow = 6
